# GetDataFromExcel issue fixed for AddAsset Excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Contacts" to "AddAsset"
$ws.Name = "AddAsset"

# Fix existing product name and append two new asset rows
$ws.Range("E2").Value = "Test Product 1"

$ws.Range("A3").Value = "AAA-5030543"
$ws.Range("B3").Value = "MA-7383994"
$ws.Range("C3").Value = "Stapler"
$ws.Range("D3").Value = "C3"
$ws.Range("E3").Value = "Test Product 2"

$ws.Range("A4").Value = "AAA-5030544"
$ws.Range("B4").Value = "MA-7383995"
$ws.Range("C4").Value = "Large Spoon"
$ws.Range("D4").Value = "C2"
$ws.Range("E4").Value = "Test Product 3"

# Widen columns C and E to fit the new data
$ws.Columns.Item(3).ColumnWidth = 12.5714285714286
$ws.Columns.Item(5).ColumnWidth = 14.4285714285714
